# Automatic update of files.
# - Bumps the "Förändrad" (changed) date in column C for every data row
#   (rows 2-20) from 2023-09-14 (serial 45183) to 2023-09-15 (serial 45184).
# - Rewrites the link-formulas in columns S:Y for the rows that carry them
#   (rows 2-4) so each HYPERLINK() call also supplies the visible link text
#   ("<Beteckning>") as its second argument - mirroring what the upstream
#   generator produced, including its malformed S-column formula and its
#   previously broken Y-column formula (which used to be stored as plain
#   inline text with a ";" argument separator instead of a real formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the "changed" date column (C) for every data row.
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}

# 2) Rewrite the hyperlink-formula columns (S:Y) for the rows that have them.
for ($row = 2; $row -le 4; $row++) {
    $art = $ws.Cells.Item($row, 1).Value2

    # S: artfynd link - upstream wrote this one malformed (unbalanced quotes / missing comma)
    $ws.Range("S$row").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/artfynd/' + $art + '.xlsx, "' + $art + '"")'

    # T: kartor link
    $ws.Range("T$row").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/kartor/' + $art + '.png", "' + $art + '")'

    # U: knärot link (row 4 never had this column - leave it alone)
    if ($row -ne 4) {
        $ws.Range("U$row").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/knärot/' + $art + '.png", "' + $art + '")'
    }

    # V: klagomål link
    $ws.Range("V$row").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomål/' + $art + '.docx", "' + $art + '")'

    # W: klagomålsmail link
    $ws.Range("W$row").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/klagomålsmail/' + $art + '.docx", "' + $art + '")'

    # X: tillsyn link
    $ws.Range("X$row").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsyn/' + $art + '.docx", "' + $art + '")'

    # Y: tillsynsmail link - previously stored as broken inline text; now a real formula
    $ws.Range("Y$row").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ANGE/tillsynsmail/' + $art + '.docx", "' + $art + '")'
}
